# Update NATMI TPM data: rows 2-7 now reflect recalculated values
# for the Sending/Target cluster combinations (ECs, FAPs, MuSCs), and
# the rows for Target cluster "ECs" (formerly rows 2-4, plus old rows
# 8-10) have been dropped from the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,20
$arr[0,0]="ECs"; $arr[0,1]="Ntf3"; $arr[0,2]="Ntrk1"; $arr[0,3]="FAPs"; $arr[0,4]=3; $arr[0,5]=1; $arr[0,6]=10.934894; $arr[0,7]=32.804682; $arr[0,8]=0.3698068269583527; $arr[0,9]=0.3698068269583527; $arr[0,10]=1; $arr[0,11]=0.3333333333333333; $arr[0,12]=0.08949833333333333; $arr[0,13]=0.268495; $arr[0,14]=0.8034178059852001; $arr[0,15]=0.8034178059852; $arr[0,16]=0.9786547881766666; $arr[0,17]=8.80789309359; $arr[0,18]=0.2971093895532283; $arr[0,19]=0.2971093895532282
$ws.Range("A2:T2").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0]="ECs"; $arr[0,1]="Ntf3"; $arr[0,2]="Ntrk1"; $arr[0,3]="MuSCs"; $arr[0,4]=3; $arr[0,5]=1; $arr[0,6]=10.934894; $arr[0,7]=32.804682; $arr[0,8]=0.3698068269583527; $arr[0,9]=0.3698068269583527; $arr[0,10]=1; $arr[0,11]=0.3333333333333333; $arr[0,12]=0.02189866666666667; $arr[0,13]=0.065696; $arr[0,14]=0.1965821940147999; $arr[0,15]=0.1965821940147999; $arr[0,16]=0.2394595987413333; $arr[0,17]=2.155136388672; $arr[0,18]=0.07269743740512444; $arr[0,19]=0.07269743740512444
$ws.Range("A3:T3").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0]="FAPs"; $arr[0,1]="Ntf3"; $arr[0,2]="Ntrk1"; $arr[0,3]="FAPs"; $arr[0,4]=3; $arr[0,5]=1; $arr[0,6]=11.451921; $arr[0,7]=34.355763; $arr[0,8]=0.3872921463699351; $arr[0,9]=0.3872921463699351; $arr[0,10]=1; $arr[0,11]=0.3333333333333333; $arr[0,12]=0.08949833333333333; $arr[0,13]=0.268495; $arr[0,14]=0.8034178059852001; $arr[0,15]=0.8034178059852; $arr[0,16]=1.024927842965; $arr[0,17]=9.224350586685; $arr[0,18]=0.3111574065118323; $arr[0,19]=0.3111574065118322
$ws.Range("A4:T4").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0]="FAPs"; $arr[0,1]="Ntf3"; $arr[0,2]="Ntrk1"; $arr[0,3]="MuSCs"; $arr[0,4]=3; $arr[0,5]=1; $arr[0,6]=11.451921; $arr[0,7]=34.355763; $arr[0,8]=0.3872921463699351; $arr[0,9]=0.3872921463699351; $arr[0,10]=1; $arr[0,11]=0.3333333333333333; $arr[0,12]=0.02189866666666667; $arr[0,13]=0.065696; $arr[0,14]=0.1965821940147999; $arr[0,15]=0.1965821940147999; $arr[0,16]=0.250781800672; $arr[0,17]=2.257036206048; $arr[0,18]=0.07613473985810289; $arr[0,19]=0.07613473985810289
$ws.Range("A5:T5").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0]="MuSCs"; $arr[0,1]="Ntf3"; $arr[0,2]="Ntrk1"; $arr[0,3]="FAPs"; $arr[0,4]=3; $arr[0,5]=1; $arr[0,6]=7.182390333333333; $arr[0,7]=21.547171; $arr[0,8]=0.2429010266717122; $arr[0,9]=0.2429010266717122; $arr[0,10]=1; $arr[0,11]=0.3333333333333333; $arr[0,12]=0.08949833333333333; $arr[0,13]=0.268495; $arr[0,14]=0.8034178059852001; $arr[0,15]=0.8034178059852; $arr[0,16]=0.6428119641827778; $arr[0,17]=5.785307677644999; $arr[0,18]=0.1951510099201395; $arr[0,19]=0.1951510099201395
$ws.Range("A6:T6").Value = $arr

$arr = New-Object 'object[,]' 1,20
$arr[0,0]="MuSCs"; $arr[0,1]="Ntf3"; $arr[0,2]="Ntrk1"; $arr[0,3]="MuSCs"; $arr[0,4]=3; $arr[0,5]=1; $arr[0,6]=7.182390333333333; $arr[0,7]=21.547171; $arr[0,8]=0.2429010266717122; $arr[0,9]=0.2429010266717122; $arr[0,10]=1; $arr[0,11]=0.3333333333333333; $arr[0,12]=0.02189866666666667; $arr[0,13]=0.065696; $arr[0,14]=0.1965821940147999; $arr[0,15]=0.1965821940147999; $arr[0,16]=0.1572847717795555; $arr[0,17]=2.155136388672; $arr[0,18]=0.04775001675157261; $arr[0,19]=0.04775001675157261
$ws.Range("A7:T7").Value = $arr


# Old rows 8-10 (Target cluster "ECs" for Sending cluster "MuSCs", plus
# stray duplicate rows) are no longer part of the table.
$ws.Range("A8:T10").EntireRow.Delete()
